# Configuration update after IOS 2015
#
# Applies the changes described by the commit:
#   - icf.185      : BO1:BU1-style "+BO2-693960" date-header formulas already
#                     existed; the equivalent block on icf_a.185 is re-entered
#                     as a single range formula so it collapses into a shared
#                     formula group (same as Excel does when you fill the same
#                     formula across adjacent cells).
#   - icf_a.185    : BU9 0.34 -> 0.342 ; BU12 1580 -> 1570
#   - Eventos.185  : A66:A71 "+B66-693960" style formulas re-entered the same
#                     way so they collapse into a shared formula group.
#   - selections on icf.185 / icf_a.185 moved to reflect where the user was
#     last working (BT27 and BU12 respectively).

$wb = $excel.ActiveWorkbook

$wsIcf    = $wb.Worksheets.Item("icf.185")
$wsIcfA   = $wb.Worksheets.Item("icf_a.185")
$wsEvent  = $wb.Worksheets.Item("Eventos.185")

# --- icf_a.185: collapse BO1:BU1 into a shared formula (values unchanged) ---
$wsIcfA.Range("BO1:BU1").Formula = "=+BO2-693960"

# --- icf_a.185: updated figures ---
$wsIcfA.Range("BU9").Value = 0.34200000000000003
$wsIcfA.Range("BU12").Value = 1570

# --- Eventos.185: collapse A66:A71 into a shared formula (values unchanged) ---
$wsEvent.Range("A66:A71").Formula = "=+B66-693960"

# --- selections: leave the cursor where the author left it ---
# Touching icf.185 activates it momentarily; re-selecting on icf_a.185
# afterwards restores it as the active/tabSelected sheet, matching the
# workbook's unchanged activeTab.
$wsIcf.Range("BT27").Select()
$wsIcfA.Range("BU12").Select()
